# Aggiornamento fino a 21 marzo
# Appends 4 new daily rows (230-233) to Sheet1, following the same pattern
# (date serial in column A styled like the previous rows, plain numbers in B:D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$newRows = @(
    @(44304, 1, 19, 452.7043126042411),
    @(44305, 0, 5, 119.1327138432213),
    @(44306, 0, 3, 71.47962830593281),
    @(44307, 0, 3, 71.47962830593281)
)

$lastRow = 229
$destRow = $lastRow + 1

foreach ($row in $newRows) {
    # Copy the formatting (number format, font, borders, alignment) of the
    # last existing date cell so the new date cells match the existing style.
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$destRow").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($destRow, 1).Value2 = $row[0]
    $ws.Cells.Item($destRow, 2).Value2 = $row[1]
    $ws.Cells.Item($destRow, 3).Value2 = $row[2]
    $ws.Cells.Item($destRow, 4).Value2 = $row[3]

    $destRow = $destRow + 1
}
